# Apply the "term 2.0.0" update:
#  - Metadata: Version 1.1.0 -> 2.0.0, Date -> 2024-06-04T14:59:10+02:00,
#    Contact -> "Kommunernes Landsforening (http://kl.dk)"
#  - Two new worksheets ("Include from FSIII 3" / "Include from FSIII 4"),
#    cloned from the existing "Include from FSIII 2" filter table, each
#    pointing at a different concept code (G1 / G2 resp. - new uuids).

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet updates -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# ---- Clone the "Include from FSIII 2" filter sheet twice -------------------
$template = $wb.Worksheets.Item("Include from FSIII 2")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$template.Copy($null, $lastSheet)
$sheet3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3.Name = "Include from FSIII 3"
$sheet3.Range("C2").Value = "993d8f7b-fbed-4a78-90d9-6efbfa835114"

$sheet3.Copy($null, $sheet3)
$sheet4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4.Name = "Include from FSIII 4"
$sheet4.Range("C2").Value = "ff47f955-3179-446f-b211-dc29de9456e3"

# ---- Restore the originally-active tab (Metadata) --------------------------
$meta.Activate()
